$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130, shifting existing rows 130:208 down to 131:209
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new daily record
$ws.Cells.Item(130, 1).Value = 5
$ws.Cells.Item(130, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(130, 3).Value = "Maule"
$ws.Cells.Item(130, 4).Value = "2022-02-14"
$ws.Cells.Item(130, 5).Value = 7
$ws.Cells.Item(130, 6).Value = "Fruta"
$ws.Cells.Item(130, 7).Value = 100108
$ws.Cells.Item(130, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(130, 9).Value = 100108005
$ws.Cells.Item(130, 10).Value = "Piña"
$ws.Cells.Item(130, 11).Value = "Caramelo"
$ws.Cells.Item(130, 12).Value = "Segunda"
$ws.Cells.Item(130, 13).Value = 250
$ws.Cells.Item(130, 14).Value = 17000
$ws.Cells.Item(130, 15).Value = 17000
$ws.Cells.Item(130, 16).Value = 17000
$ws.Cells.Item(130, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(130, 18).Value = "Ecuador"
$ws.Cells.Item(130, 19).Value = 1214
$ws.Cells.Item(130, 20).Value = 14
